$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet1: LUT borders shift by one (more precise binning boarders) ---

# N5 had its own unique formula (not part of the shared group below):
# J5&", "  ->  J5-1&", "
$ws.Range("N5").Formula = '=IF(J5<>"",J5-1&", "&IF($H5<>$H6,CHAR(10),""),256^$C$8-1&CHAR(10))'

# N6:N26 is the shared-formula block (si=7 in the original file):
# J#&", "  ->  J#-1&", "
$ws.Range("N6:N26").Formula = '=IF(J6<>"",J6-1&", "&IF($H6<>$H7,CHAR(10),""),256^$C$8-1&CHAR(10))'

# N27 previously had its own broken formula (referencing #REF!); it now follows
# the same J-1 pattern as the rest of the column (with J27/H27 blank this still
# evaluates to the "255" terminator row, matching the rest of the block).
$ws.Range("N27").Formula = '=IF(J27<>"",J27-1&", "&IF($H27<>$H28,CHAR(10),""),256^$C$8-1&CHAR(10))'

# Selection moved from P5:P43 to N31
$ws.Range("N31").Select()
